# whisker_module_BOM.xlsx - "Added datasheets, finalized final_BOM"
#
# 1. Fix the Socket (CD74HCT151E) row's datasheet/part URL in H7: it used to
#    point at a DigiKey *search* URL; point it at the real product-detail page.
# 2. Fill in row 17 of the BOM table (previously an all-zero placeholder row)
#    with a new "Protoboard" line item, formatted with Excel's built-in "Bad"
#    cell style (as if flagged for follow-up), including a hyperlink on the
#    Frys part-number cell.
# 3. Leave the selection on J18, matching the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Correct the datasheet URL text shown in H7 ------------------------
$ws.Range("H7").Value = "http://www.digikey.com/product-detail/en/CD74HCT151E/296-2139-5-ND/38312"

# --- 2. Populate row 17 with the new Protoboard line item -----------------
$ws.Range("A17").Value = "Protoboard"
$ws.Range("B17").Value = "4x5"
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = "Twin"
$ws.Range("E17").Value = "8000-45-LF"
$ws.Range("F17").Value = "Frys"
$ws.Range("G17").Value = 4986181
$ws.Range("H17").Value = "http://www.frys.com/product/4986181?source=googleps&gclid=CNPX5OnlxbQCFcxAMgodfm0AYw"
$ws.Range("I17").Value = 12.69
$ws.Range("J17").Value = 2
# K17 already contains the calculated-column formula Price*Quantity and will
# recompute automatically to 25.38.

# Apply the built-in "Bad" cell style (red text on a pink fill) to the whole
# row, then fix up the two cells that need extra number formatting:
#  - G17 keeps a left-aligned General number (it's a part/catalog number)
#  - I17 keeps the currency format used elsewhere in the Price column
$ws.Range("A17:J17").Style = "Bad"
$ws.Range("G17").HorizontalAlignment = -4131
$ws.Range("I17").NumberFormat = """$""#,##0.00_);[Red]\(""$""#,##0.00\)"

# Give the row the same custom height the author set (14.25pt)
$ws.Rows.Item(17).RowHeight = 14.25

# Hyperlink the new Frys part number cell. Passing the URL as the
# TextToDisplay value makes Excel cache it as the hyperlink's display text
# (matching the other hyperlinks in this sheet), then we restore G17's real
# value (the numeric catalog id) afterwards.
$ws.Hyperlinks.Add($ws.Range("G17"), "http://www.frys.com/product/4986181?source=googleps&gclid=CNPX5OnlxbQCFcxAMgodfm0AYw", [Type]::Missing, [Type]::Missing, "http://www.frys.com/product/4986181?source=googleps&gclid=CNPX5OnlxbQCFcxAMgodfm0AYw") | Out-Null
$ws.Range("G17").Value = 4986181

# --- 3. Restore the saved selection state ----------------------------------
$ws.Range("J18").Select() | Out-Null
